# Generate Report for Handback
# Updates timestamps (and one status value) produced by a fresh handback
# report generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the rows
# that previously read 2016-08-25 02:15:30
$wsOverview.Range("G2").Value = "2016-08-25 02:16:18"
$wsOverview.Range("G4").Value = "2016-08-25 02:16:18"

# Priority column (E) status "ht" -> "mt" (shared by zh-cn and de-de sheets)
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime (H) and
# Correspond Handback DateTime (K)
$wsZhCn.Range("H2").Value = "2016-08-25 02:16:13"
$wsZhCn.Range("H4").Value = "2016-08-25 02:16:13"

$wsZhCn.Range("K2").Value = "2016-08-25 02:16:29"
$wsZhCn.Range("K4").Value = "2016-08-25 02:16:29"

# de-de sheet: Correspond Handoff Datetime (H)
$wsDeDe.Range("H2").Value = "2016-08-25 02:16:37"
$wsDeDe.Range("H4").Value = "2016-08-25 02:16:37"
